# The deck shipped with its two theme parts swapped relative to how they
# are used: ppt/theme/theme1.xml (the slide master's theme, i.e. the design
# applied to every slide) held the "Integral" / "Red Violet" palette, while
# ppt/theme/theme2.xml (the notes master's theme) held the default "Office
# Theme" palette. The edit swaps the two 12-slot colour schemes so the
# slide master now carries the standard Office colours.
#
# PowerPoint's theme colours (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) are
# reached through Slide.ThemeColorScheme, which maps 1:1 onto the
# <a:clrScheme> children of ppt/theme/theme1.xml (the theme backing the
# slide master/layouts used by every slide). RGB is encoded the same way
# VBA's RGB() macro does: 0x00BBGGRR, i.e. R + G*256 + B*65536.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Hex-RGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office Theme" colours previously living in theme2.xml.
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $tcs.Colors($i).RGB = Hex-RGB $officeColors[$i - 1]
}
